$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (G) values - computed std/mean derived s_vals
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 2
    12 = 3
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 2
    27 = 2
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 3
    39 = 3
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 1
    45 = 2
    46 = 2
    47 = 2
    48 = 1
    49 = 2
    50 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
